$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date and fill in the Jurisdiction ---
$meta = $wb.Worksheets.Item("Metadata")

# Row 8  -> Property "Date"         / Value (B8)
# Row 11 -> Property "Jurisdiction" / Value (B11)
$meta.Range("B8").Value = "2025-07-11T12:29:53+00:00"
$meta.Range("B11").Value = "FRANCE"

$wb.Save()
